# The commit this fixture corresponds to ("Moving from POI 3.17.0 to
# 4.0.1") only re-serialised the package: every hunk in the diff swaps
# the attribute/namespace-declaration order on elements such as
# <w:document>, <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>,
# <w:latentStyles>, <w:lsdException> and <w:style> - the element names,
# attribute names and attribute values are all identical before and
# after, only the order in which the upgraded XML writer emitted the
# attributes changed. There is no reachable Word object-model
# operation (text, formatting, styles, structure, ...) that this
# corresponds to, so no edit is required here.

$d = $word.ActiveDocument
